# Scheduled market-data refresh for the Leviathan_Profits workbook.
# Updates the currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the affected Leve rows across all eight job sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) to the latest pulled market-board values.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 617.8  # H96: was 615.3333
$ws.Cells.Item(96, 9).Value = 735  # I96: was 844.1667
$ws.Cells.Item(96, 10).Value = 149  # J96: was 157.66667
$ws.Cells.Item(96, 11).Value = 2205  # K96: was 2532.5001
$ws.Cells.Item(96, 12).Value = 447  # L96: was 473.00001
$ws.Cells.Item(96, 13).Value = -832  # M96: was -1159.5001
$ws.Cells.Item(96, 14).Value = -3193  # N96: was -3219.00001
$ws.Cells.Item(100, 8).Value = 4603.647  # H100: was 3636.25
$ws.Cells.Item(100, 9).Value = 3919.16  # I100: was 3950.84
$ws.Cells.Item(100, 10).Value = 6505  # J100: was 2512.7144
$ws.Cells.Item(100, 11).Value = 3919.16  # K100: was 3950.84
$ws.Cells.Item(100, 12).Value = 6505  # L100: was 2512.7144
$ws.Cells.Item(100, 13).Value = -3378.16  # M100: was -3409.84
$ws.Cells.Item(100, 14).Value = -7587  # N100: was -3594.7144
$ws.Cells.Item(113, 8).Value = 75921.57000000001  # H113: was 81595.16
$ws.Cells.Item(113, 9).Value = 253024.75  # I113: was 203386.4
$ws.Cells.Item(113, 10).Value = 5080.3  # J113: was 5475.625
$ws.Cells.Item(113, 11).Value = 253024.75  # K113: was 203386.4
$ws.Cells.Item(113, 12).Value = 5080.3  # L113: was 5475.625
$ws.Cells.Item(113, 13).Value = -249770.75  # M113: was -200132.4
$ws.Cells.Item(113, 14).Value = -11588.3  # N113: was -11983.625
$ws.Cells.Item(132, 8).Value = 2052.5576  # H132: was 2090.255
$ws.Cells.Item(132, 9).Value = 817.7347  # I132: was 832.0625
$ws.Cells.Item(132, 11).Value = 2453.2041  # K132: was 2496.1875
$ws.Cells.Item(132, 13).Value = 76.79590000000007  # M132: was 33.8125
$ws.Cells.Item(135, 8).Value = 1569.3478  # H135: was 1570.5217
$ws.Cells.Item(135, 9).Value = 1449.375  # I135: was 1451.0625
$ws.Cells.Item(135, 11).Value = 13044.375  # K135: was 13059.5625
$ws.Cells.Item(135, 13).Value = -10509.375  # M135: was -10524.5625
$ws.Cells.Item(137, 8).Value = 4043.4443  # H137: was 3635.3635
$ws.Cells.Item(137, 9).Value = 2898.5  # I137: was 2698.4285
$ws.Cells.Item(137, 10).Value = 6333.3335  # J137: was 5275
$ws.Cells.Item(137, 11).Value = 8695.5  # K137: was 8095.2855
$ws.Cells.Item(137, 12).Value = 19000.0005  # L137: was 15825
$ws.Cells.Item(137, 13).Value = -6145.5  # M137: was -5545.2855
$ws.Cells.Item(137, 14).Value = -24100.0005  # N137: was -20925
$ws.Cells.Item(139, 8).Value = 179944  # H139: was 179998.67
$ws.Cells.Item(139, 10).Value = 179944  # J139: was 179998.67
$ws.Cells.Item(139, 12).Value = 179944  # L139: was 179998.67
$ws.Cells.Item(139, 14).Value = -190224  # N139: was -190278.67
$ws.Cells.Item(141, 8).Value = 1321.5  # H141: was 2097
$ws.Cells.Item(141, 9).Value = 1367.7142  # I141: was 2097
$ws.Cells.Item(141, 10).Value = 998  # J141: was 0
$ws.Cells.Item(141, 11).Value = 4103.142599999999  # K141: was 6291
$ws.Cells.Item(141, 12).Value = 2994  # L141: was 0
$ws.Cells.Item(141, 13).Value = 1076.857400000001  # M141: was -1111
$ws.Cells.Item(141, 14).Value = -13354  # N141: was None

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1295.6538  # H2: was 1295.6154
$ws.Cells.Item(2, 9).Value = 1236.05  # I2: was 1236
$ws.Cells.Item(2, 11).Value = 1236.05  # K2: was 1236
$ws.Cells.Item(2, 13).Value = -1123.05  # M2: was -1123
$ws.Cells.Item(32, 8).Value = 4922.1787  # H32: was 5099.8477
$ws.Cells.Item(32, 9).Value = 3921.7158  # I32: was 4121.4165
$ws.Cells.Item(32, 10).Value = 17499.428  # J32: was 15373.375
$ws.Cells.Item(32, 11).Value = 3921.7158  # K32: was 4121.4165
$ws.Cells.Item(32, 12).Value = 17499.428  # L32: was 15373.375
$ws.Cells.Item(32, 13).Value = -3634.7158  # M32: was -3834.4165
$ws.Cells.Item(32, 14).Value = -18073.428  # N32: was -15947.375
$ws.Cells.Item(45, 8).Value = 8018.3687  # H45: was 8813.941000000001
$ws.Cells.Item(45, 9).Value = 9586.923000000001  # I45: was 11101.637
$ws.Cells.Item(45, 11).Value = 9586.923000000001  # K45: was 11101.637
$ws.Cells.Item(45, 13).Value = -9209.923000000001  # M45: was -10724.637
$ws.Cells.Item(68, 8).Value = 25099  # H68: was 24399.334
$ws.Cells.Item(68, 10).Value = 25099  # J68: was 24399.334
$ws.Cells.Item(68, 12).Value = 25099  # L68: was 24399.334
$ws.Cells.Item(68, 14).Value = -26721  # N68: was -26021.334
$ws.Cells.Item(71, 8).Value = 25099  # H71: was 24399.334
$ws.Cells.Item(71, 10).Value = 25099  # J71: was 24399.334
$ws.Cells.Item(71, 12).Value = 75297  # L71: was 73198.00199999999
$ws.Cells.Item(71, 14).Value = -83409  # N71: was -81310.00199999999
$ws.Cells.Item(74, 8).Value = 1657  # H74: was 1624.6072
$ws.Cells.Item(74, 9).Value = 948.0526  # I74: was 938.15
$ws.Cells.Item(74, 11).Value = 948.0526  # K74: was 938.15
$ws.Cells.Item(74, 13).Value = -74.05259999999998  # M74: was -64.14999999999998
$ws.Cells.Item(77, 8).Value = 1657  # H77: was 1624.6072
$ws.Cells.Item(77, 9).Value = 948.0526  # I77: was 938.15
$ws.Cells.Item(77, 11).Value = 4740.263  # K77: was 4690.75
$ws.Cells.Item(77, 13).Value = -372.2629999999999  # M77: was -322.75
$ws.Cells.Item(110, 8).Value = 2332.0789  # H110: was 2305.6667
$ws.Cells.Item(110, 9).Value = 1350.6897  # I110: was 1349.0667
$ws.Cells.Item(110, 11).Value = 1350.6897  # K110: was 1349.0667
$ws.Cells.Item(110, 13).Value = 694.3103000000001  # M110: was 695.9332999999999
$ws.Cells.Item(116, 8).Value = 1295.6538  # H116: was 1295.6154
$ws.Cells.Item(116, 9).Value = 1236.05  # I116: was 1236
$ws.Cells.Item(116, 11).Value = 1236.05  # K116: was 1236
$ws.Cells.Item(116, 13).Value = 1057.95  # M116: was 1058

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1295.6538  # H3: was 1295.6154
$ws.Cells.Item(3, 9).Value = 1236.05  # I3: was 1236
$ws.Cells.Item(3, 11).Value = 1236.05  # K3: was 1236
$ws.Cells.Item(3, 13).Value = -1122.05  # M3: was -1122
$ws.Cells.Item(24, 8).Value = 743.5  # H24: was 381
$ws.Cells.Item(24, 9).Value = 743.5  # I24: was 381
$ws.Cells.Item(24, 11).Value = 743.5  # K24: was 381
$ws.Cells.Item(24, 13).Value = -508.5  # M24: was -146
$ws.Cells.Item(80, 8).Value = 2516.6875  # H80: was 2671.8
$ws.Cells.Item(80, 9).Value = 598.25  # I80: was 734.3333
$ws.Cells.Item(80, 11).Value = 598.25  # K80: was 734.3333
$ws.Cells.Item(80, 13).Value = 399.75  # M80: was 263.6667
$ws.Cells.Item(83, 8).Value = 2516.6875  # H83: was 2671.8
$ws.Cells.Item(83, 9).Value = 598.25  # I83: was 734.3333
$ws.Cells.Item(83, 11).Value = 2991.25  # K83: was 3671.6665
$ws.Cells.Item(83, 13).Value = 2000.75  # M83: was 1320.3335
$ws.Cells.Item(99, 8).Value = 2423.4075  # H99: was 2186.1936
$ws.Cells.Item(99, 9).Value = 2457.32  # I99: was 2199.0688
$ws.Cells.Item(99, 11).Value = 2457.32  # K99: was 2199.0688
$ws.Cells.Item(99, 13).Value = -959.3200000000002  # M99: was -701.0688
$ws.Cells.Item(134, 8).Value = 94843.69  # H134: was 108301.29
$ws.Cells.Item(134, 9).Value = 104521.34  # I134: was 121142.28
$ws.Cells.Item(134, 11).Value = 313564.02  # K134: was 363426.84
$ws.Cells.Item(134, 13).Value = -311029.02  # M134: was -360891.84

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1154.1  # H16: was 1236.8889
$ws.Cells.Item(16, 9).Value = 807.25  # I16: was 940
$ws.Cells.Item(16, 11).Value = 807.25  # K16: was 940
$ws.Cells.Item(16, 13).Value = -520.25  # M16: was -653
$ws.Cells.Item(62, 8).Value = 0  # H62: was 3000
$ws.Cells.Item(62, 10).Value = 0  # J62: was 3000
$ws.Cells.Item(62, 12).ClearContents()  # L62: was 3000, now blank
$ws.Cells.Item(62, 14).Value = 0  # N62: was -4248
$ws.Cells.Item(65, 8).Value = 0  # H65: was 3000
$ws.Cells.Item(65, 10).Value = 0  # J65: was 3000
$ws.Cells.Item(65, 12).ClearContents()  # L65: was 15000, now blank
$ws.Cells.Item(65, 14).Value = 0  # N65: was -21240
$ws.Cells.Item(107, 8).Value = 2381.7896  # H107: was 2193.6667
$ws.Cells.Item(107, 9).Value = 1804  # I107: was 1687
$ws.Cells.Item(107, 10).Value = 3176.25  # J107: was 2869.2222
$ws.Cells.Item(107, 11).Value = 1804  # K107: was 1687
$ws.Cells.Item(107, 12).Value = 3176.25  # L107: was 2869.2222
$ws.Cells.Item(107, 13).Value = 116  # M107: was 233
$ws.Cells.Item(107, 14).Value = -7016.25  # N107: was -6709.2222
$ws.Cells.Item(113, 8).Value = 1154.1  # H113: was 1236.8889
$ws.Cells.Item(113, 9).Value = 807.25  # I113: was 940
$ws.Cells.Item(113, 11).Value = 807.25  # K113: was 940
$ws.Cells.Item(113, 13).Value = 1362.75  # M113: was 1230
$ws.Cells.Item(132, 8).Value = 3175.9119  # H132: was 3176
$ws.Cells.Item(132, 9).Value = 3066.1  # I132: was 3066.2
$ws.Cells.Item(132, 11).Value = 9198.299999999999  # K132: was 9198.599999999999
$ws.Cells.Item(132, 13).Value = -6668.299999999999  # M132: was -6668.599999999999
$ws.Cells.Item(141, 8).Value = 251454.22  # H141: was 237857.2
$ws.Cells.Item(141, 10).Value = 251454.22  # J141: was 237857.2
$ws.Cells.Item(141, 12).Value = 251454.22  # L141: was 237857.2
$ws.Cells.Item(141, 14).Value = -261814.22  # N141: was -248217.2

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 43.75  # H7: was 437.91666
$ws.Cells.Item(7, 9).Value = 36.666668  # I7: was 471.81818
$ws.Cells.Item(7, 11).Value = 110.000004  # K7: was 1415.45454
$ws.Cells.Item(7, 13).Value = 1.999995999999996  # M7: was -1303.45454
$ws.Cells.Item(22, 8).Value = 2618.8  # H22: was 3390
$ws.Cells.Item(22, 10).Value = 2618.8  # J22: was 3390
$ws.Cells.Item(22, 12).Value = 7856.400000000001  # L22: was 10170
$ws.Cells.Item(22, 14).Value = -8194.400000000001  # N22: was -10508
$ws.Cells.Item(27, 8).Value = 2618.8  # H27: was 3390
$ws.Cells.Item(27, 10).Value = 2618.8  # J27: was 3390
$ws.Cells.Item(27, 12).Value = 7856.400000000001  # L27: was 10170
$ws.Cells.Item(27, 14).Value = -8060.400000000001  # N27: was -10374
$ws.Cells.Item(64, 8).Value = 2386.25  # H64: was 2210
$ws.Cells.Item(64, 10).Value = 4000  # J64: was 3360
$ws.Cells.Item(64, 12).Value = 12000  # L64: was 10080
$ws.Cells.Item(64, 14).Value = -12540  # N64: was -10620
$ws.Cells.Item(67, 8).Value = 2386.25  # H67: was 2210
$ws.Cells.Item(67, 10).Value = 4000  # J67: was 3360
$ws.Cells.Item(67, 12).Value = 12000  # L67: was 10080
$ws.Cells.Item(67, 14).Value = -13872  # N67: was -11952
$ws.Cells.Item(80, 8).Value = 4891.364  # H80: was 4838.6924
$ws.Cells.Item(80, 10).Value = 4880.5  # J80: was 4825.25
$ws.Cells.Item(80, 12).Value = 14641.5  # L80: was 14475.75
$ws.Cells.Item(80, 14).Value = -16513.5  # N80: was -16347.75
$ws.Cells.Item(83, 8).Value = 4891.364  # H83: was 4838.6924
$ws.Cells.Item(83, 10).Value = 4880.5  # J83: was 4825.25
$ws.Cells.Item(83, 12).Value = 43924.5  # L83: was 43427.25
$ws.Cells.Item(83, 14).Value = -53284.5  # N83: was -52787.25
$ws.Cells.Item(94, 8).Value = 501562  # H94: was 203756.8
$ws.Cells.Item(94, 9).Value = 501562  # I94: was 335982
$ws.Cells.Item(94, 10).Value = 0  # J94: was 5419
$ws.Cells.Item(94, 11).Value = 1504686  # K94: was 1007946
$ws.Cells.Item(94, 12).Value = 0  # L94: was 16257
$ws.Cells.Item(94, 13).ClearContents()  # M94: was -1007270, now blank
$ws.Cells.Item(94, 14).Value = -1504010  # N94: was -17609

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 5000199.5  # H3: was 3750216
$ws.Cells.Item(3, 9).Value = 5000049.5  # I3: was 4000045.5
$ws.Cells.Item(3, 10).Value = 5000500  # J3: was 3333833.2
$ws.Cells.Item(3, 11).Value = 5000049.5  # K3: was 4000045.5
$ws.Cells.Item(3, 12).Value = 5000500  # L3: was 3333833.2
$ws.Cells.Item(3, 13).Value = -4999933.5  # M3: was -3999929.5
$ws.Cells.Item(3, 14).Value = -5000732  # N3: was -3334065.2
$ws.Cells.Item(113, 8).Value = 2428.625  # H113: was 2302.6667
$ws.Cells.Item(113, 9).Value = 2203.6  # I113: was 2221.8
$ws.Cells.Item(113, 10).Value = 2803.6667  # J113: was 2403.75
$ws.Cells.Item(113, 11).Value = 2203.6  # K113: was 2221.8
$ws.Cells.Item(113, 12).Value = 2803.6667  # L113: was 2403.75
$ws.Cells.Item(113, 13).Value = -33.59999999999991  # M113: was -51.80000000000018
$ws.Cells.Item(113, 14).Value = -7143.6667  # N113: was -6743.75
$ws.Cells.Item(132, 8).Value = 4880.0884  # H132: was 4997.6665
$ws.Cells.Item(132, 9).Value = 3007.652  # I132: was 3098.9092
$ws.Cells.Item(132, 11).Value = 9022.956  # K132: was 9296.7276
$ws.Cells.Item(132, 13).Value = -6492.956  # M132: was -6766.7276

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 19427.424  # H46: was 18686.777
$ws.Cells.Item(46, 9).Value = 27801  # I46: was 24936.63
$ws.Cells.Item(46, 10).Value = 3610.6667  # J46: was 3843.375
$ws.Cells.Item(46, 11).Value = 27801  # K46: was 24936.63
$ws.Cells.Item(46, 12).Value = 3610.6667  # L46: was 3843.375
$ws.Cells.Item(46, 13).Value = -27613  # M46: was -24748.63
$ws.Cells.Item(46, 14).Value = -3986.6667  # N46: was -4219.375
$ws.Cells.Item(140, 8).Value = 89252.336  # H140: was 80498
$ws.Cells.Item(140, 10).Value = 89252.336  # J140: was 80498
$ws.Cells.Item(140, 12).Value = 89252.336  # L140: was 80498
$ws.Cells.Item(140, 14).Value = -99612.336  # N140: was -90858

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3956.5  # H62: was 3991.7144
$ws.Cells.Item(62, 10).Value = 3809  # J62: was 3858.5
$ws.Cells.Item(62, 12).Value = 3809  # L62: was 3858.5
$ws.Cells.Item(62, 14).Value = -5057  # N62: was -5106.5
$ws.Cells.Item(65, 8).Value = 3956.5  # H65: was 3991.7144
$ws.Cells.Item(65, 10).Value = 3809  # J65: was 3858.5
$ws.Cells.Item(65, 12).Value = 19045  # L65: was 19292.5
$ws.Cells.Item(65, 14).Value = -25285  # N65: was -25532.5
$ws.Cells.Item(107, 8).Value = 38463304  # H107: was 41668556
$ws.Cells.Item(107, 9).Value = 2334.8572  # I107: was 2677
$ws.Cells.Item(107, 11).Value = 7004.571599999999  # K107: was 8031
$ws.Cells.Item(107, 13).Value = -5084.571599999999  # M107: was -6111
$ws.Cells.Item(113, 8).Value = 533.5  # H113: was 420.83334
$ws.Cells.Item(113, 9).Value = 542  # I113: was 415
$ws.Cells.Item(113, 10).Value = 525  # J113: was 450
$ws.Cells.Item(113, 11).Value = 1626  # K113: was 1245
$ws.Cells.Item(113, 12).Value = 1575  # L113: was 1350
$ws.Cells.Item(113, 13).Value = 544  # M113: was 925
$ws.Cells.Item(113, 14).Value = -5915  # N113: was -5690
$ws.Cells.Item(132, 8).Value = 3754.889  # H132: was 3755.5557
$ws.Cells.Item(132, 9).Value = 2994.3044  # I132: was 2995.348
$ws.Cells.Item(132, 11).Value = 8982.913199999999  # K132: was 8986.044
$ws.Cells.Item(132, 13).Value = -6452.913199999999  # M132: was -6456.044
